$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44425
$ws.Range("J2").Value = 30

$ws.Range("D3").Value = 44676
$ws.Range("J3").Value = 40
$ws.Range("K3").Value = 30000
$ws.Range("M3").Value = 30000
$ws.Range("O3").Value = 'Provincia de Quillota'
$ws.Range("P3").Value = 3000

$ws.Range("D4").Value = 44446
$ws.Range("J4").Value = 40
$ws.Range("K4").Value = 27000
$ws.Range("L4").Value = 27000
$ws.Range("M4").Value = 27000
$ws.Range("P4").Value = 2700

$ws.Range("D5").Value = 44454
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 25000
$ws.Range("L5").Value = 25000
$ws.Range("M5").Value = 25000
$ws.Range("P5").Value = 2500

$ws.Range("D6").Value = 44694
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 30000
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = 30000
$ws.Range("P6").Value = 3000

$ws.Range("D7").Value = 44410
$ws.Range("J7").Value = 50

$ws.Range("D8").Value = 44405
$ws.Range("J8").Value = 40

$ws.Range("D9").Value = 44355
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 23000
$ws.Range("L9").Value = 24000
$ws.Range("M9").Value = 23400
$ws.Range("P9").Value = 2340

$ws.Range("D10").Value = 44392
$ws.Range("J10").Value = 25

$ws.Range("D11").Value = 44417
$ws.Range("J11").Value = 15
$ws.Range("K11").Value = 25000
$ws.Range("L11").Value = 25000
$ws.Range("M11").Value = 25000
$ws.Range("O11").Value = 'Provincia de Quillota'
$ws.Range("P11").Value = 2500

$ws.Range("D12").Value = 44349
$ws.Range("J12").Value = 45

$ws.Range("D13").Value = 44354
$ws.Range("J13").Value = 30
$ws.Range("K13").Value = 24000
$ws.Range("L13").Value = 24000
$ws.Range("M13").Value = 24000
$ws.Range("P13").Value = 2400

$ws.Range("D14").Value = 44677
$ws.Range("J14").Value = 20

$ws.Range("D15").Value = 44365
$ws.Range("J15").Value = 85
$ws.Range("K15").Value = 22000
$ws.Range("L15").Value = 22000
$ws.Range("M15").Value = 22000
$ws.Range("P15").Value = 2200

$ws.Range("D16").Value = 44715
$ws.Range("J16").Value = 12
$ws.Range("K16").Value = 28000
$ws.Range("L16").Value = 28000
$ws.Range("M16").Value = 28000
$ws.Range("P16").Value = 2800

$ws.Range("D17").Value = 44719
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = 28000
$ws.Range("L17").Value = 30000
$ws.Range("M17").Value = 29333
$ws.Range("P17").Value = 2933

$ws.Range("D18").Value = 44721
$ws.Range("J18").Value = 100
$ws.Range("O18").Value = 'Región Metropolitana'

$ws.Range("D19").Value = 44462
$ws.Range("J19").Value = 50
$ws.Range("K19").Value = 25000
$ws.Range("L19").Value = 25000
$ws.Range("M19").Value = 25000
$ws.Range("O19").Value = 'Provincia de Quillota'
$ws.Range("P19").Value = 2500

$ws.Range("D20").Value = 44665
$ws.Range("J20").Value = 10

$ws.Range("D21").Value = 44412
$ws.Range("J21").Value = 50

$ws.Range("D22").Value = 44461
$ws.Range("J22").Value = 40
$ws.Range("K22").Value = 25000
$ws.Range("L22").Value = 25000
$ws.Range("M22").Value = 25000
$ws.Range("P22").Value = 2500

$ws.Range("D23").Value = 44467
$ws.Range("J23").Value = 40
$ws.Range("K23").Value = 23000
$ws.Range("M23").Value = 24000
$ws.Range("P23").Value = 2400

$ws.Range("D24").Value = 44669
$ws.Range("J24").Value = 15
$ws.Range("K24").Value = 30000
$ws.Range("L24").Value = 30000
$ws.Range("M24").Value = 30000
$ws.Range("O24").Value = 'Región Metropolitana'
$ws.Range("P24").Value = 3000

$ws.Range("D25").Value = 44664
$ws.Range("J25").Value = 20

$ws.Range("D26").Value = 44420
$ws.Range("J26").Value = 55

$ws.Range("D27").Value = 44376
$ws.Range("J27").Value = 45
$ws.Range("K27").Value = 23000
$ws.Range("L27").Value = 23000
$ws.Range("M27").Value = 23000
$ws.Range("P27").Value = 2300

$ws.Range("D28").Value = 44397
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = 27000
$ws.Range("L28").Value = 27000
$ws.Range("M28").Value = 27000
$ws.Range("O28").Value = 'Provincia de Quillota'
$ws.Range("P28").Value = 2700

$ws.Range("D29").Value = 44428
$ws.Range("J29").Value = 30

$ws.Range("D30").Value = 44431
$ws.Range("J30").Value = 65

$ws.Range("D31").Value = 44701
$ws.Range("K31").Value = 28000
$ws.Range("L31").Value = 30000
$ws.Range("M31").Value = 29333
$ws.Range("O31").Value = 'Región Metropolitana'
$ws.Range("P31").Value = 2933

$ws.Range("D32").Value = 44406
$ws.Range("J32").Value = 120
$ws.Range("K32").Value = 24000
$ws.Range("M32").Value = 24542
$ws.Range("P32").Value = 2454

$ws.Range("D33").Value = 44438
$ws.Range("J33").Value = 50

$ws.Range("D34").Value = 44466
$ws.Range("J34").Value = 50
$ws.Range("K34").Value = 25000
$ws.Range("L34").Value = 25000
$ws.Range("M34").Value = 25000
$ws.Range("P34").Value = 2500

$ws.Range("D35").Value = 44704
$ws.Range("J35").Value = 30
$ws.Range("K35").Value = 30000
$ws.Range("L35").Value = 30000
$ws.Range("M35").Value = 30000
$ws.Range("O35").Value = 'Región Metropolitana'
$ws.Range("P35").Value = 3000

$ws.Range("D36").Value = 44413
$ws.Range("J36").Value = 40

$ws.Range("D37").Value = 44729
$ws.Range("J37").Value = 45
$ws.Range("K37").Value = 28000
$ws.Range("L37").Value = 28000
$ws.Range("M37").Value = 28000
$ws.Range("P37").Value = 2800

$ws.Range("D38").Value = 44427
$ws.Range("J38").Value = 40
$ws.Range("K38").Value = 25000
$ws.Range("L38").Value = 25000
$ws.Range("M38").Value = 25000
$ws.Range("P38").Value = 2500

$ws.Range("D39").Value = 44468

$ws.Range("D40").Value = 44434
$ws.Range("J40").Value = 55
$ws.Range("K40").Value = 25000
$ws.Range("L40").Value = 25000
$ws.Range("M40").Value = 25000
$ws.Range("P40").Value = 2500

$ws.Range("D41").Value = 44448
$ws.Range("J41").Value = 15
$ws.Range("L41").Value = 25000
$ws.Range("M41").Value = 25000
$ws.Range("P41").Value = 2500

$ws.Range("D42").Value = 44441
$ws.Range("J42").Value = 70

$ws.Range("D43").Value = 44727
$ws.Range("J43").Value = 95
$ws.Range("K43").Value = 28000
$ws.Range("L43").Value = 28000
$ws.Range("M43").Value = 28000
$ws.Range("P43").Value = 2800

$ws.Range("D44").Value = 44708
$ws.Range("J44").Value = 20
$ws.Range("K44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("M44").Value = 30000
$ws.Range("O44").Value = 'Región Metropolitana'
$ws.Range("P44").Value = 3000

$ws.Range("D45").Value = 44419
$ws.Range("J45").Value = 25
$ws.Range("K45").Value = 25000
$ws.Range("L45").Value = 25000
$ws.Range("M45").Value = 25000
$ws.Range("P45").Value = 2500

$ws.Range("D46").Value = 44447
$ws.Range("J46").Value = 30
$ws.Range("K46").Value = 27000
$ws.Range("L46").Value = 27000
$ws.Range("M46").Value = 27000
$ws.Range("P46").Value = 2700

$ws.Range("D47").Value = 44396
$ws.Range("J47").Value = 20

$ws.Range("D48").Value = 44400
$ws.Range("J48").Value = 12
$ws.Range("K48").Value = 24000
$ws.Range("L48").Value = 24000
$ws.Range("M48").Value = 24000
$ws.Range("P48").Value = 2400

$ws.Range("D49").Value = 44670
$ws.Range("J49").Value = 25
$ws.Range("K49").Value = 30000
$ws.Range("L49").Value = 30000
$ws.Range("M49").Value = 30000
$ws.Range("O49").Value = 'Región Metropolitana'
$ws.Range("P49").Value = 3000

$ws.Range("D50").Value = 44722
$ws.Range("J50").Value = 20
$ws.Range("K50").Value = 28000
$ws.Range("L50").Value = 28000
$ws.Range("M50").Value = 28000
$ws.Range("O50").Value = 'Región Metropolitana'
$ws.Range("P50").Value = 2800

$ws.Range("D51").Value = 44453
$ws.Range("J51").Value = 40
$ws.Range("K51").Value = 27000
$ws.Range("L51").Value = 27000
$ws.Range("M51").Value = 27000
$ws.Range("P51").Value = 2700

$ws.Range("D52").Value = 44487

$ws.Range("D53").Value = 44432
$ws.Range("J53").Value = 15

$ws.Range("D54").Value = 44371
$ws.Range("J54").Value = 50

$ws.Range("D55").Value = 44691
$ws.Range("J55").Value = 55
$ws.Range("K55").Value = 30000
$ws.Range("L55").Value = 30000
$ws.Range("M55").Value = 30000
$ws.Range("P55").Value = 3000

$ws.Range("D56").Value = 44389
$ws.Range("J56").Value = 65
$ws.Range("O56").Value = 'Provincia de Quillota'

$ws.Range("D57").Value = 44463
$ws.Range("J57").Value = 40
$ws.Range("K57").Value = 25000
$ws.Range("M57").Value = 26000
$ws.Range("P57").Value = 2600

$ws.Range("D58").Value = 44685
$ws.Range("J58").Value = 20
$ws.Range("K58").Value = 30000
$ws.Range("L58").Value = 30000
$ws.Range("M58").Value = 30000
$ws.Range("P58").Value = 3000

$ws.Range("D59").Value = 44449
$ws.Range("J59").Value = 12

$ws.Range("D60").Value = 44455
$ws.Range("J60").Value = 20

$ws.Range("D61").Value = 44379
$ws.Range("J61").Value = 35
$ws.Range("K61").Value = 22000
$ws.Range("L61").Value = 22000
$ws.Range("M61").Value = 22000
$ws.Range("O61").Value = 'Provincia de Quillota'
$ws.Range("P61").Value = 2200

$ws.Range("D62").Value = 44356
$ws.Range("J62").Value = 15
$ws.Range("K62").Value = 24000
$ws.Range("L62").Value = 24000
$ws.Range("M62").Value = 24000
$ws.Range("P62").Value = 2400

$ws.Range("D63").Value = 44720
$ws.Range("J63").Value = 30
$ws.Range("K63").Value = 28000
$ws.Range("L63").Value = 28000
$ws.Range("M63").Value = 28000
$ws.Range("O63").Value = 'Región Metropolitana'
$ws.Range("P63").Value = 2800

$ws.Range("D64").Value = 44350
$ws.Range("J64").Value = 40
$ws.Range("K64").Value = 24000
$ws.Range("M64").Value = 24375
$ws.Range("P64").Value = 2438

$ws.Range("D65").Value = 44452
$ws.Range("J65").Value = 80
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = 25000
$ws.Range("P65").Value = 2500

$ws.Range("D66").Value = 44382
$ws.Range("J66").Value = 50
$ws.Range("K66").Value = 25000
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = 25000
$ws.Range("P66").Value = 2500

$ws.Range("D67").Value = 44469
$ws.Range("J67").Value = 80
$ws.Range("K67").Value = 25000
$ws.Range("L67").Value = 25000
$ws.Range("M67").Value = 25000
$ws.Range("O67").Value = 'Provincia de Quillota'
$ws.Range("P67").Value = 2500

$ws.Range("D68").Value = 44435
$ws.Range("J68").Value = 185
$ws.Range("K68").Value = 25000
$ws.Range("L68").Value = 27000
$ws.Range("M68").Value = 25162
$ws.Range("P68").Value = 2516

$ws.Range("D69").Value = 44390
$ws.Range("J69").Value = 15
$ws.Range("K69").Value = 25000
$ws.Range("L69").Value = 25000
$ws.Range("M69").Value = 25000
$ws.Range("P69").Value = 2500

$ws.Range("D70").Value = 44372
$ws.Range("J70").Value = 20
$ws.Range("K70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("M70").Value = 25000
$ws.Range("P70").Value = 2500

$ws.Range("D71").Value = 44384
$ws.Range("J71").Value = 40
$ws.Range("K71").Value = 25000
$ws.Range("M71").Value = 25000
$ws.Range("P71").Value = 2500

$ws.Range("D72").Value = 44706
$ws.Range("J72").Value = 50
$ws.Range("K72").Value = 30000
$ws.Range("L72").Value = 30000
$ws.Range("M72").Value = 30000
$ws.Range("O72").Value = 'Región Metropolitana'
$ws.Range("P72").Value = 3000

$ws.Range("D73").Value = 44411
$ws.Range("J73").Value = 40

$ws.Range("D74").Value = 44476
$ws.Range("J74").Value = 50
$ws.Range("K74").Value = 25000
$ws.Range("L74").Value = 25000
$ws.Range("M74").Value = 25000
$ws.Range("P74").Value = 2500

$ws.Range("D75").Value = 44474
$ws.Range("J75").Value = 30
$ws.Range("O75").Value = 'Región Metropolitana'

$ws.Range("D76").Value = 44726
$ws.Range("J76").Value = 35
$ws.Range("O76").Value = 'Provincia de Quillota'

$ws.Range("D77").Value = 44426
$ws.Range("J77").Value = 30
$ws.Range("K77").Value = 25000
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = 25000
$ws.Range("O77").Value = 'Provincia de Quillota'
$ws.Range("P77").Value = 2500

$ws.Range("D78").Value = 44477
$ws.Range("J78").Value = 40

$ws.Range("D79").Value = 44473
$ws.Range("J79").Value = 50

$ws.Range("D80").Value = 44699
$ws.Range("J80").Value = 95
$ws.Range("K80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("M80").Value = 30000
$ws.Range("P80").Value = 3000

$ws.Range("D81").Value = 44385
$ws.Range("J81").Value = 80
$ws.Range("K81").Value = 25000
$ws.Range("L81").Value = 25000
$ws.Range("M81").Value = 25000
$ws.Range("O81").Value = 'Provincia de Quillota'
$ws.Range("P81").Value = 2500

$ws.Range("D82").Value = 44348
$ws.Range("J82").Value = 3

$ws.Range("D83").Value = 44386
$ws.Range("J83").Value = 20

$ws.Range("D84").Value = 44690
$ws.Range("J84").Value = 40
$ws.Range("K84").Value = 30000
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = 30000
$ws.Range("P84").Value = 3000

$ws.Range("D85").Value = 44433
$ws.Range("J85").Value = 25
$ws.Range("K85").Value = 25000
$ws.Range("M85").Value = 25000
$ws.Range("P85").Value = 2500

$ws.Range("D86").Value = 44707
$ws.Range("J86").Value = 80
$ws.Range("K86").Value = 30000
$ws.Range("L86").Value = 30000
$ws.Range("M86").Value = 30000
$ws.Range("P86").Value = 3000

$ws.Range("D87").Value = 44421
$ws.Range("J87").Value = 55

$ws.Range("D88").Value = 44483
$ws.Range("J88").Value = 20
$ws.Range("K88").Value = 25000
$ws.Range("L88").Value = 25000
$ws.Range("M88").Value = 25000
$ws.Range("P88").Value = 2500
